# Groups0.xlsx edit: grow the carpool children list from 5 to 9 kids and
# refresh their relative-position coordinates (the data source the
# map/graph colouring reads from - see commit message about matching the
# "school" point colour across the graphs), shifting the trailing
# school/cost/time summary rows down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used as padding after some first/last names in the
# source data (matches the rest of the sheet's existing name fields).
$nbsp = [char]0x00A0

function Set-CellText($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    if ($value -match '^-?\d+(\.\d+)?$') {
        # Force text so numeric-looking tokens ("30.0", "9", ...) keep their
        # exact literal text instead of being coerced to a Number by Excel.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}

function Set-ChildRow($row, $values) {
    for ($col = 1; $col -le $values.Count; $col++) {
        Set-CellText $row $col $values[$col - 1]
    }
}

# nChildren: 5 -> 9
Set-CellText 4 2 "9"

# Children table (rows 6-14), 9 children (id 0..8)
Set-ChildRow 6  @("0", "12", ("Frankie " + $nbsp), ("Flavia " + $nbsp), "-0.18,-7.37", "Cyrus(mother): 0522363358", "7:00:00", "32.0")
Set-ChildRow 7  @("1", "20", "Ron", "Cohen", "-2.33,-7.05", "Bernardine(mother): 0576270618", "7:03:00", "29.0")
Set-ChildRow 8  @("2", "3", ("Alexia " + $nbsp), ("Ramonita " + $nbsp), "-5.15,-3.61", "Han(father): 0567537032", "7:08:00", "24.0")
Set-ChildRow 9  @("3", "11", ("Randolph " + $nbsp), ("Bridgette " + $nbsp), "-9.16,-4.56", "Lenny(father): 0505536740", "7:13:00", "19.0")
Set-ChildRow 10 @("4", "17", ("Britta " + $nbsp), ("Jamel " + $nbsp), "-6.72,-1.52", "Albertine(father): 0574981040", "7:18:00", "14.0")
Set-ChildRow 11 @("5", "4", ("Francisca " + $nbsp), ("Stevie " + $nbsp), "-7.14,-1.26", "Bernardine(mother): 0561339273", "7:19:00", "13.0")
Set-ChildRow 12 @("6", "10", ("Demetra " + $nbsp), ("Francene " + $nbsp), "-6.3,-0.62", "Dorian(mother): 0534328089", "7:21:00", "11.0")
Set-ChildRow 13 @("7", "16", ("Collette " + $nbsp), ("Billi " + $nbsp), "-4.83,-1.02", "Elias(mother): 0578741979", "7:24:00", "8.0")
Set-ChildRow 14 @("8", "5", ("Patti " + $nbsp), ("Lavenia " + $nbsp), "-0.63,-1.53", "Jennell(mother): 0503029941", "7:29:00", "3.0")

# school row moves from 11 -> 15; only the pickup time (G) changes.
Set-ChildRow 15 @("school", "3", "Ironiah", "mySchool", "0,0", "Shir(secretary): 0523345098", "7:32:00")

# cost row moves from 12 -> 16 (values unchanged).
Set-ChildRow 16 @("cost", "25")

# time row moves from 13 -> 17; B changes 30.0 -> 32.0.
Set-ChildRow 17 @("time", "32.0")
